# Add a new "total" sheet to the database workbook, after the existing
# "cost_log" sheet (i.e. at the end), and give it the header row used by
# the new update_total_table() method: id, month_id, year, save, cost.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "total"

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "month_id"
$ws.Range("C1").Value = "year"
$ws.Range("D1").Value = "save"
$ws.Range("E1").Value = "cost"

$ws.Range("F9").Select() | Out-Null
